# Repull data, push all data, mean calculation
# Update the "dSF" column (F) values for the rows where the recalculated
# mean-based delta differs from 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3"  = -2
    "F5"  = 1
    "F8"  = -2
    "F9"  = -3
    "F11" = 3
    "F12" = -3
    "F14" = -1
    "F15" = 4
    "F16" = -2
    "F17" = -3
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
